$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row at row 9 (pushes FEDESAFOROL/URAID-N/ZYRTEC + the
#    totals row + the footer row down by one), then populate it with the new
#    "ENEMAX ENEMA 120 ML" item.
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# Copy the formatting of an existing, fully-styled data row into the blank
# row so the new row's cell styles match the rest of the table.
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("9:9").RowHeight = 25.5

# Fill in the new row's values.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "ENEMAX ENEMA 120 ML"
$ws.Range("H9").Value = "5:0"
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = "40.00"
$ws.Range("P9").Value = "40.0000"
$ws.Range("Q9").Value = "1:0"

# Renumber the "م" counter column for the rows that followed (they shifted
# down by one row).
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# ---------------------------------------------------------------------------
# 2. Update the running total (old row 12, now row 13) to reflect the new
#    item's price.
# ---------------------------------------------------------------------------
$ws.Range("P13").Value = 220.90000000000001

# ---------------------------------------------------------------------------
# 3. Bump the printed timestamp by a minute.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Saturday, 13 September, 2025 10:15 AM"

# ---------------------------------------------------------------------------
# 4. Rebuild the merged-cell ranges so the new row's merges exist and the
#    merge list is in the same (row-major) order as the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Cells.UnMerge()

$ws.Range("D2:N2").Merge()
$ws.Range("A3:Q3").Merge()
$ws.Range("F4:H4").Merge()
$ws.Range("I4:K4").Merge()
$ws.Range("L4:N4").Merge()
$ws.Range("B6:G6").Merge()
$ws.Range("H6:L6").Merge()
$ws.Range("N6:O6").Merge()
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()
$ws.Range("P13:Q13").Merge()
$ws.Range("A14:F14").Merge()
$ws.Range("G14:I14").Merge()
$ws.Range("K14:Q14").Merge()

Write-Host "Edit complete"
